# Re-balanced point cost, new character
#
# 1) The deck's "last updated" date field (Insert > Header & Footer >
#    Date and time > Update automatically) is re-stamped on the slide
#    master and on every custom layout that carries its own copy of the
#    placeholder.
# 2) On the character sheet itself, the DODGE stat (first cell of the
#    small DODGE/ARMOR/SAVE table) is re-balanced from 4 to 6.

$p = $ppt.ActivePresentation

$oldDate = "28.11.2016"
$newDate = "22.01.2017"

# --- Slide master date placeholder -----------------------------------
$master = $p.SlideMaster
foreach ($shp in $master.Shapes) {
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq $oldDate) {
            $tr.Text = $newDate
        }
    }
}

# --- Every custom layout's date placeholder ---------------------------
foreach ($layout in $master.CustomLayouts) {
    foreach ($shp in $layout.Shapes) {
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# --- Character sheet: bump the DODGE stat from 4 to 6 ------------------
$slide = $p.Slides.Item(1)
foreach ($shp in $slide.Shapes) {
    if ($shp.HasTable -and $shp.Name -eq "Table 13") {
        $cell = $shp.Table.Cell(1, 1)
        if ($cell.Shape.TextFrame.TextRange.Text -eq "4") {
            $cell.Shape.TextFrame.TextRange.Text = "6"
        }
    }
}
